# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row number => new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 1260
    3  = 1155
    4  = 884
    7  = 644
    11 = 2291
    12 = 1574
    13 = 1307
    15 = 228
    17 = 738
    18 = 28
    19 = 276
    22 = 8
    24 = 4456
    25 = 208
    26 = 17
    31 = 6
    32 = 645
    37 = 363
    39 = 124
    41 = 123
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - row number => new F value
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 1260
    5  = 1155
    6  = 884
    11 = 644
    17 = 2291
    18 = 1574
    19 = 1307
    21 = 228
    24 = 738
    25 = 28
    26 = 276
    28 = 8
    29 = 4456
    30 = 208
    31 = 17
    36 = 6
    37 = 645
    41 = 363
    43 = 124
    45 = 123
}
foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
